# Thêm dữ liệu mẫu cho Template_Import_Student.xlsx
# (mirrors the "thêm dữ liệu mẫu cho file excel và fix đường dẫn bình luận" commit)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hàng 2: dữ liệu mẫu tương ứng với các cột ở hàng 1
# A = Họ và tên, B = Mã số sinh viên, C = Email, D = Ngày sinh, E = Giới tính
$ws.Range("A2").Value = "Nguyễn Văn A"
$ws.Range("B2").Value = "0306200000"
$ws.Range("C2").Value = "0306200000@caothang.edu.vn"
$ws.Range("D2").Value = 32874
$ws.Range("E2").Value = "Nam"

# Ghi chú dữ liệu mẫu ở cột G
$ws.Range("G2").Value = "(Đây là dữ liệu mẫu, vui lòng hãy xóa để thêm mới)"

# Gắn hyperlink (mailto) cho ô email mẫu
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:0306200000@caothang.edu.vn") | Out-Null

# Cập nhật vùng chọn hiện hành giống bản đã lưu
$ws.Range("J8").Select() | Out-Null
